# Attendance SPRING workbook update:
# - Rename "Liam Nestelroad" -> "Megan McGinns" in row 10 (B10) on each sheet
# - Delete the now-duplicate "Megan McGinns" row (row 11) on each sheet
# - Record additional meeting attendance (checkmarks) and month/day headers
# - Adjust data validation ranges to the new (shrunk) data extents
# - Restore per-sheet selection / active cell

$wb = $excel.ActiveWorkbook
$check = [char]0x2714

# ---------------------------------------------------------------------------
# SPONSOR sheet
# ---------------------------------------------------------------------------
$sponsor = $wb.Worksheets.Item("SPONSOR")

$sponsor.Range("B10").Value = "Megan McGinns"
$sponsor.Rows.Item(11).Delete() | Out-Null

$sponsor.Range("D3").Value = "Jan"
$sponsor.Range("D4").Value = 26
$sponsor.Range("D5:D10").Value = $check

$sponsor.Range("C5:N11").Validation.Delete() | Out-Null
$sponsor.Range("I5:N10").Validation.Add(3, 1, 1, "SPONSOR!checkbox") | Out-Null
$sponsor.Range("C5:H10").Validation.Add(3, 1, 1, "SPONSOR!checkbox") | Out-Null

# ---------------------------------------------------------------------------
# TEAM sheet
# ---------------------------------------------------------------------------
$team = $wb.Worksheets.Item("TEAM")

$team.Range("B10").Value = "Megan McGinns"
$team.Rows.Item(11).Delete() | Out-Null

$team.Range("C3").Value = "Jan"
$team.Range("D3").Value = "Jan"
$team.Range("C4").Value = 19
$team.Range("D4").Value = 26
$team.Range("C5:C10").Value = $check

$team.Range("C5:N11").Validation.Delete() | Out-Null
$team.Range("D5:N10").Validation.Add(3, 1, 1, "TEAM!checkbox") | Out-Null
$team.Range("C5:C10").Validation.Add(3, 1, 1, "TEAM!checkbox") | Out-Null

# ---------------------------------------------------------------------------
# TA sheet
# ---------------------------------------------------------------------------
$ta = $wb.Worksheets.Item("TA")

$ta.Range("B10").Value = "Megan McGinns"
$ta.Rows.Item(11).Delete() | Out-Null

$ta.Range("C3").Value = "Jan"
$ta.Range("C4").Value = 26
$ta.Range("C5:C8").Value = $check
$ta.Range("C10").Value = $check

$ta.Range("C5:N11").Validation.Delete() | Out-Null
$ta.Range("D5:N10").Validation.Add(3, 1, 1, "TEAM!checkbox") | Out-Null
$ta.Range("C5:C10").Validation.Add(3, 1, 1, "SPONSOR!checkbox") | Out-Null

# ---------------------------------------------------------------------------
# Selections (match author's last on-screen state; TA stays the active sheet)
# ---------------------------------------------------------------------------
$sponsor.Activate() | Out-Null
$sponsor.Range("D5").Select() | Out-Null

$team.Activate() | Out-Null
$team.Range("D4").Select() | Out-Null

$ta.Activate() | Out-Null
$ta.Range("B19").Select() | Out-Null
